$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 106194537
$ws.Range("B3").Value = 56411
$ws.Range("D3").Value = 'NT'
$ws.Range("E3").Value = 100049
$ws.Range("F3").Value = 'Spillkråka'
$ws.Range("G3").Value = 'Dryocopus martius'
$ws.Range("H3").Value = '(Linnaeus, 1758)'
$ws.Range("M3").Value = 'äldre spår'
$ws.Range("Q3").Value = 540471.4570663463
$ws.Range("R3").Value = 7198931.455751203
$ws.Range("Z3").Value = '11:05'
$ws.Range("AB3").Value = '11:05'
$ws.Range("A4").Value = 106194518
$ws.Range("B4").Value = 78596
$ws.Range("D4").Value = 'LC'
$ws.Range("E4").Value = 6462
$ws.Range("F4").Value = 'Stuplav'
$ws.Range("G4").Value = 'Nephroma bellum'
$ws.Range("H4").Value = '(Spreng.) Tuck.'
$ws.Range("Q4").Value = 540220.9713190208
$ws.Range("R4").Value = 7199288.915839214
$ws.Range("Z4").Value = '12:25'
$ws.Range("AB4").Value = '12:25'
$ws.Range("A5").Value = 106194523
$ws.Range("B5").Value = 77506
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 6425
$ws.Range("F5").Value = 'Garnlav'
$ws.Range("G5").Value = 'Alectoria sarmentosa'
$ws.Range("H5").Value = '(Ach.) Ach.'
$ws.Range("Q5").Value = 540266.2807869902
$ws.Range("R5").Value = 7199332.822708054
$ws.Range("Z5").Value = '12:21'
$ws.Range("AB5").Value = '12:21'
$ws.Range("A6").Value = 106194526
$ws.Range("B6").Value = 78596
$ws.Range("D6").Value = 'LC'
$ws.Range("E6").Value = 6462
$ws.Range("F6").Value = 'Stuplav'
$ws.Range("G6").Value = 'Nephroma bellum'
$ws.Range("H6").Value = '(Spreng.) Tuck.'
$ws.Range("Q6").Value = 540458.3903157733
$ws.Range("R6").Value = 7199332.018118081
$ws.Range("Z6").Value = '11:50'
$ws.Range("AB6").Value = '11:50'
$ws.Range("AJ6").ClearContents()
$ws.Range("AK6").ClearContents()
$ws.Range("AO6").ClearContents()
$ws.Range("AQ6").ClearContents()
$ws.Range("AR6").ClearContents()
$ws.Range("A7").Value = 106194569
$ws.Range("B7").Value = 76487
$ws.Range("D7").Value = 'VU'
$ws.Range("E7").Value = 1794
$ws.Range("F7").Value = 'Rödskaftad svartspik'
$ws.Range("G7").Value = 'Chaenothecopsis haematopus'
$ws.Range("H7").Value = 'Tibell'
$ws.Range("I7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("Q7").Value = 540458.4132686888
$ws.Range("R7").Value = 7199330.320494538
$ws.Range("Z7").Value = '11:49'
$ws.Range("AB7").Value = '11:49'
$ws.Range("AJ7").Value = 'sälg'
$ws.Range("AK7").Value = 'Salix caprea'
$ws.Range("AO7").Value = 'På sälgved och sälgticka # Salix caprea'
$ws.Range("AQ7").Value = 'Isak Vahlström'
$ws.Range("AR7").Value = "'2301261149"
$ws.Range("A8").Value = 106194540
$ws.Range("B8").Value = 73693
$ws.Range("E8").Value = 6440
$ws.Range("F8").Value = 'Vitgrynig nållav'
$ws.Range("G8").Value = 'Chaenotheca subroscida'
$ws.Range("H8").Value = '(Eitner) Zahlbr.'
$ws.Range("Q8").Value = 540471.7153629591
$ws.Range("R8").Value = 7198912.354880349
$ws.Range("Z8").Value = '11:02'
$ws.Range("AB8").Value = '11:02'
$ws.Range("A9").Value = 106194528
$ws.Range("B9").Value = 55608
$ws.Range("E9").Value = 102612
$ws.Range("F9").Value = 'Järpe'
$ws.Range("G9").Value = 'Tetrastes bonasia'
$ws.Range("H9").Value = '(Linnaeus, 1758)'
$ws.Range("I9").Value = "'1"
$ws.Range("M9").Value = 'spel/sång'
$ws.Range("Q9").Value = 540421.9978498913
$ws.Range("R9").Value = 7199131.588012774
$ws.Range("Z9").Value = '11:27'
$ws.Range("AB9").Value = '11:27'
$ws.Range("A10").Value = 106194536
$ws.Range("B10").Value = 76489
$ws.Range("D10").Value = 'DD'
$ws.Range("E10").Value = 6000248
$ws.Range("F10").Value = 'Mörk kådsvartspik'
$ws.Range("G10").Value = 'Chaenothecopsis montana'
$ws.Range("H10").Value = 'Rikkinen'
$ws.Range("Q10").Value = 540443.4200716898
$ws.Range("R10").Value = 7198930.227664446
$ws.Range("Z10").Value = '11:11'
$ws.Range("AB10").Value = '11:11'
$ws.Range("A11").Value = 106194541
$ws.Range("Q11").Value = 540467.4483503525
$ws.Range("R11").Value = 7198913.570809621
$ws.Range("Z11").Value = '11:02'
$ws.Range("AB11").Value = '11:02'
$ws.Range("A12").Value = 106194529
$ws.Range("B12").Value = 78570
$ws.Range("D12").Value = 'NT'
$ws.Range("E12").Value = 2081
$ws.Range("F12").Value = 'Skrovellav'
$ws.Range("G12").Value = 'Lobaria scrobiculata'
$ws.Range("H12").Value = '(Scop.) DC.'
$ws.Range("Q12").Value = 540414.6648215668
$ws.Range("R12").Value = 7199139.554680186
$ws.Range("Z12").Value = '11:26'
$ws.Range("AB12").Value = '11:26'
$ws.Range("A13").Value = 106194522
$ws.Range("B13").Value = 81236
$ws.Range("E13").Value = 1312
$ws.Range("F13").Value = 'Gammelgransskål'
$ws.Range("G13").Value = 'Pseudographis pinicola'
$ws.Range("H13").Value = '(Nyl.) Rehm'
$ws.Range("Q13").Value = 540266.2979199347
$ws.Range("R13").Value = 7199331.549491068
$ws.Range("A14").Value = 106194520
$ws.Range("B14").Value = 76504
$ws.Range("E14").Value = 314
$ws.Range("F14").Value = 'Vitskaftad svartspik'
$ws.Range("G14").Value = 'Chaenothecopsis viridialba'
$ws.Range("H14").Value = '(Kremp.) A.F.W.Schmidt'
$ws.Range("Q14").Value = 540266.2693650301
$ws.Range("R14").Value = 7199333.671519079
$ws.Range("Z14").Value = '12:21'
$ws.Range("AB14").Value = '12:21'
$ws.Range("A15").Value = 106194517
$ws.Range("B15").Value = 77506
$ws.Range("E15").Value = 6425
$ws.Range("F15").Value = 'Garnlav'
$ws.Range("G15").Value = 'Alectoria sarmentosa'
$ws.Range("H15").Value = '(Ach.) Ach.'
$ws.Range("M15").ClearContents()
$ws.Range("Q15").Value = 540220.9713190208
$ws.Range("R15").Value = 7199288.915839214
$ws.Range("Z15").Value = '12:25'
$ws.Range("AB15").Value = '12:25'
$ws.Range("A16").Value = 106194519
$ws.Range("B16").Value = 78603
$ws.Range("D16").Value = 'LC'
$ws.Range("E16").Value = 6464
$ws.Range("F16").Value = 'Luddlav'
$ws.Range("G16").Value = 'Nephroma resupinatum'
$ws.Range("H16").Value = '(L.) Ach.'
$ws.Range("Q16").Value = 540217.9683172886
$ws.Range("R16").Value = 7199290.997918672
$ws.Range("Z16").Value = '12:25'
$ws.Range("AB16").Value = '12:25'
$ws.Range("A17").Value = 106194521
$ws.Range("B17").Value = 73693
$ws.Range("D17").Value = 'NT'
$ws.Range("E17").Value = 6440
$ws.Range("F17").Value = 'Vitgrynig nållav'
$ws.Range("G17").Value = 'Chaenotheca subroscida'
$ws.Range("H17").Value = '(Eitner) Zahlbr.'
$ws.Range("Q17").Value = 540266.2807869902
$ws.Range("R17").Value = 7199332.822708054
$ws.Range("Z17").Value = '12:21'
$ws.Range("AB17").Value = '12:21'
$ws.Range("A18").Value = 106194539
$ws.Range("B18").Value = 76490
$ws.Range("E18").Value = 228579
$ws.Range("F18").Value = 'Liten svartspik'
$ws.Range("G18").Value = 'Chaenothecopsis nana'
$ws.Range("H18").Value = 'Tibell'
$ws.Range("Q18").Value = 540470.4519041032
$ws.Range("R18").Value = 7198911.488709519
$ws.Range("A19").Value = 106194534
$ws.Range("B19").Value = 77506
$ws.Range("E19").Value = 6425
$ws.Range("F19").Value = 'Garnlav'
$ws.Range("G19").Value = 'Alectoria sarmentosa'
$ws.Range("H19").Value = '(Ach.) Ach.'
$ws.Range("Q19").Value = 540443.3971279609
$ws.Range("R19").Value = 7198931.925514392
$ws.Range("Z19").Value = '11:12'
$ws.Range("AB19").Value = '11:12'
$ws.Range("A20").Value = 106194538
$ws.Range("B20").Value = 73507
$ws.Range("D20").Value = 'LC'
$ws.Range("E20").Value = 6428
$ws.Range("F20").Value = 'Rostfläck'
$ws.Range("G20").Value = 'Arthonia vinosa'
$ws.Range("H20").Value = 'Leight.'
$ws.Range("Q20").Value = 540471.5259453177
$ws.Range("R20").Value = 7198926.362197533
$ws.Range("Z20").Value = '11:04'
$ws.Range("AB20").Value = '11:04'
$ws.Range("A21").Value = 106194527
$ws.Range("B21").Value = 77506
$ws.Range("D21").Value = 'NT'
$ws.Range("E21").Value = 6425
$ws.Range("F21").Value = 'Garnlav'
$ws.Range("G21").Value = 'Alectoria sarmentosa'
$ws.Range("H21").Value = '(Ach.) Ach.'
$ws.Range("Q21").Value = 540458.3903157733
$ws.Range("R21").Value = 7199332.018118081
$ws.Range("Z21").Value = '11:50'
$ws.Range("AB21").Value = '11:50'
